# Replace Sheet1 with a new Sheet2 (mirrors the original authoring workflow,
# which bumps sheetId the same way a delete+add of a sheet would in Excel).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$old = $wb.Worksheets.Item("Sheet1")
$old.Delete()

# ---- Content -----------------------------------------------------------
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "What the Science Says"

$ws.Range("A2").Value = "`"Climate's changed before`""
$ws.Range("B2").Value = "Climate reacts to whatever forces it to change at the time; humans are now the dominant forcing."

$ws.Range("A3").Value = "`"It's the sun`""
$ws.Range("B3").Value = "In the last 35 years of global warming, sun and climate have been going in opposite directions"

$ws.Range("A4").Value = "`"It's not bad`""
$ws.Range("B4").Value = "Negative impacts of global warming on agriculture, health & environment far outweigh any positives."

# ---- Column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24.8
$ws.Columns.Item(2).ColumnWidth = 92.1

# ---- Styling --------------------------------------------------------------
# Theme colours used below (iron_native's 1-based ThemeColor index -> OOXML
# theme attribute): 2 -> theme"0" (Background 1/white), 10 -> theme"9"
# (Accent 6, the green used throughout this theme).
$headerGreen = 4697456      # 70AD47 (Accent 6)
$lightGreen  = 14348258     # E2EFDA (Accent 6, Lighter 80%)
$borderGreen = 9359785      # A9D18E (Accent 6, Lighter 40%)

# Row 1: bold white header text on a solid green fill.
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Font.ThemeColor = 2
$headerRange.Interior.ThemeColor = 10

# Row 2 and row 4: light-green banded rows.
$bandRange = $ws.Range("A2:B2,A4:B4")
$bandRange.Interior.Color = $lightGreen

# Apply the thin green borders (top+bottom on every data cell, plus a right
# edge on column B to close off the little "table").
foreach ($rowNum in 1..4) {
    $left = $ws.Cells.Item($rowNum, 1)
    $right = $ws.Cells.Item($rowNum, 2)

    $left.Borders.Item(8).LineStyle = 1
    $left.Borders.Item(8).Weight = 2
    $left.Borders.Item(8).Color = $borderGreen
    $left.Borders.Item(9).LineStyle = 1
    $left.Borders.Item(9).Weight = 2
    $left.Borders.Item(9).Color = $borderGreen

    $right.Borders.Item(8).LineStyle = 1
    $right.Borders.Item(8).Weight = 2
    $right.Borders.Item(8).Color = $borderGreen
    $right.Borders.Item(9).LineStyle = 1
    $right.Borders.Item(9).Weight = 2
    $right.Borders.Item(9).Color = $borderGreen
    $right.Borders.Item(10).LineStyle = 1
    $right.Borders.Item(10).Weight = 2
    $right.Borders.Item(10).Color = $borderGreen
}

# ---- Selection -------------------------------------------------------------
$ws.Range("B12").Select()

Write-Host "done"
